$d = $word.ActiveDocument

function Replace-Exact($old, $new) {
    # Anchored, case-sensitive, whole-document Find/Replace. Safe to use when the
    # run being edited is the last run in its paragraph (or is followed only by
    # runs with different formatting), so no unrelated run gets coalesced into it.
    $range = $d.Content
    $ok = $range.Find.Execute($old, $true, $false, $false, $false, $false, `
                               $true, 1, $false, $new, 2)
    if (-not $ok) {
        throw "Find.Execute failed to find: $old"
    }
}

function Find-Range($text) {
    # Returns a Range positioned exactly over the first occurrence of $text
    # (case-sensitive), without modifying anything.
    $range = $d.Content
    $ok = $range.Find.Execute($text, $true, $false, $false, $false, $false, `
                               $true, 1, $false)
    if (-not $ok) {
        throw "Find.Execute failed to find: $text"
    }
    return $range
}

# --- Paragraph "Reconhecimento e reconhecimento / de marca limitados: ..." ---
Replace-Exact "Reconhecimento e reconhecimento" "Reconhecimento e conscientização da marca limitados"
Replace-Exact " de marca limitados: Alcançar visibilidade nesses novos mercados é um obstáculo primário, exigindo esforços de marketing robustos para construir a presença da marca Adatum desde o início." ": alcançar visibilidade nesses novos mercados é um grande obstáculo, exigindo esforços robustos de marketing para construir a presença de marca da Adatum do zero."

# --- Paragraph "Concorrência / intensa: ..." ---
# This paragraph has two more runs after the one being edited (a single space,
# then "A Adatum deve articular..."), sharing the exact same run formatting as the
# run we are about to edit. A plain Find/Replace (or any Range.Text write) on
# that run causes the engine to coalesce those identically-formatted trailing
# runs into it, which would incorrectly alter parts of the document outside the
# scope of this edit.
#
# To stop that: give the trailing runs a (temporarily) *different* format right
# before the text edit (so they are not candidates for merging while the edit
# happens), then restore their original format afterwards as its own edit. That
# leaves them as independent <w:r> elements, byte-identical to the originals.
$oldSentence = " intensa: O setor de serviços em nuvem no Canadá é ferozmente competitivo, com vários players."
$newSentence = ": o setor de serviços em nuvem no Canadá é ferozmente competitivo, com vários envolvidos."

$restOfSentence = Find-Range $oldSentence
$tailStart = $restOfSentence.End
$spaceRun = $d.Range($tailStart, $tailStart + 1)
$afterSpace = Find-Range "A Adatum deve articular claramente o valor único de suas soluções para conquistar espaço no mercado."

$spaceRun.Font.Italic = $true
$afterSpace.Font.Italic = $true

$restOfSentence2 = $d.Range($restOfSentence.Start, $restOfSentence.End)
$restOfSentence2.Text = $newSentence

$newTailStart = $restOfSentence2.Start + $newSentence.Length
$spaceRun2 = $d.Range($newTailStart, $newTailStart + 1)
$afterSpace2 = $d.Range($newTailStart + 1, $newTailStart + 1 + 101)
$spaceRun2.Font.Italic = $false
$afterSpace2.Font.Italic = $false

Replace-Exact "Concorrência" "Concorrência intensa"

# --- Paragraph "Preferências e expectativas diversificadas / dos clientes: ..." ---
Replace-Exact "Preferências e expectativas diversificadas" "Preferências e expectativas diversas do cliente"
Replace-Exact " dos clientes: Adaptar os produtos e o marketing para se alinhar às variadas demandas desses mercados é crucial para ressoar com as empresas e consumidores locais." ": adaptar produtos e marketing para se alinhar às diversas demandas desses mercados é crucial para identificação com as empresas e os consumidores locais."

# --- Paragraph "Desafios / regulatórios e de conformidade: ..." ---
Replace-Exact "Desafios" "Desafios de regulamentação e conformidade"
Replace-Exact " regulatórios e de conformidade: A Adatum enfrenta a complexa tarefa de navegar pelas distintas regulamentações de privacidade, segurança e operação de dados da região, exigindo esforços diligentes de conformidade." ": a Adatum enfrenta a complexa tarefa de navegar pelos distintos regulamentos operacionais e de privacidade de dados da região, o que demanda esforços diligentes de conformidade."

# --- Paragraph "Complexidades operacionais e logísticas / : O estabelecimento ..." ---
Replace-Exact ": O estabelecimento de operações eficientes e inter-regionais apresenta desafios logísticos, especialmente na manutenção de altos níveis de serviço e no gerenciamento de data centers em localizações geográficas." ": estabelecer operações inter-regionais eficientes apresenta desafios logísticos, especialmente na manutenção de altos níveis de serviço e no gerenciamento de data centers em localizações geográficas."

Write-Output "done"
